$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F14").Value = "93_referral_statement"
$ws.Range("F22").Value = "18_hazards_to_humans_and_domestic_animals"
$ws.Range("F24").Value = "ppe"
$ws.Range("F27").Value = "off target movement || application instructions || env warning - species || env warning - water"
$ws.Range("F31").Value = "135_product_information"
$ws.Range("F34").Value = "use restrictions"
$ws.Range("F36").Value = "135_product_information"
$ws.Range("F37").Value = "mixing"
$ws.Range("F38").Value = "mixing"
$ws.Range("F39").Value = "mixing"
$ws.Range("F51").Value = "154_pesticide_storage"
